$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value edits -------------------------------------------------
$ws.Range("A1").Value = "sdasd"
$ws.Range("B2").Value = "29.4.2020"
$ws.Range("E2").Value = "ssaa22"
$ws.Range("B3").Value = "29.4.2020"
$ws.Range("E3").Value = 0
$ws.Range("B5").Value = 22
$ws.Range("C5").Value = 22

# --- Unmerge the vertical blocks first -----------------------------------
# Clear()/ClearContents() silently no-op on a range that overlaps a merged
# block, so every merge touching A6:E13 has to be broken before the clear
# below runs. A5:A7 / B5:B7 get reduced to single-cell merges afterwards.
$ws.Range("A5:A7").UnMerge()
$ws.Range("B5:B7").UnMerge()
$ws.Range("A8:A9").UnMerge()
$ws.Range("B8:B9").UnMerge()
$ws.Range("A10:A11").UnMerge()
$ws.Range("B10:B11").UnMerge()
$ws.Range("A12:A13").UnMerge()
$ws.Range("B12:B13").UnMerge()

# --- Drop the extra product rows (6-13) ---------------------------------
# Clear() wipes content + formatting so the cells disappear entirely and
# the rows collapse back down to bare placeholders, matching rows 14-37.
$ws.Range("A6:E13").Clear()

# --- Re-merge A5/B5 as single-cell ranges (mirrors the saved mergeCells) -
$ws.Range("B5").Merge()
$ws.Range("A5").Merge()
